$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3099.8667
$ws.Range("I28").Value = 1798.25
$ws.Range("K28").Value = 1798.25
$ws.Range("M28").Value = -1313.25
$ws.Range("H39").Value = 1044.9524
$ws.Range("I39").Value = 613.5454999999999
$ws.Range("K39").Value = 1840.6365
$ws.Range("M39").Value = -1544.6365
$ws.Range("H43").Value = 557108.4
$ws.Range("I43").Value = 13726.875
$ws.Range("J43").Value = 1178115.9
$ws.Range("K43").Value = 13726.875
$ws.Range("L43").Value = 1178115.9
$ws.Range("M43").Value = -13657.875
$ws.Range("N43").Value = -1178253.9
$ws.Range("H98").Value = 1993.1372
$ws.Range("I98").Value = 2134.5854
$ws.Range("J98").Value = 1413.2
$ws.Range("K98").Value = 2134.5854
$ws.Range("L98").Value = 1413.2
$ws.Range("M98").Value = -636.5853999999999
$ws.Range("N98").Value = -4409.2
$ws.Range("H122").Value = 1993.1372
$ws.Range("I122").Value = 2134.5854
$ws.Range("J122").Value = 1413.2
$ws.Range("K122").Value = 6403.7562
$ws.Range("L122").Value = 4239.6
$ws.Range("M122").Value = -3953.7562
$ws.Range("N122").Value = -9139.6
$ws.Range("H129").Value = 1630.5555
$ws.Range("I129").Value = 636
$ws.Range("J129").Value = 2426.2
$ws.Range("K129").Value = 1908
$ws.Range("L129").Value = 7278.599999999999
$ws.Range("M129").Value = 3092
$ws.Range("N129").Value = -17278.6
$ws.Range("H132").Value = 1616.5
$ws.Range("I132").Value = 1599.1333
$ws.Range("K132").Value = 4797.3999
$ws.Range("M132").Value = -2267.3999
$ws.Range("H133").Value = 78780
$ws.Range("J133").Value = 78780
$ws.Range("L133").Value = 78780
$ws.Range("N133").Value = -88900
$ws.Range("H138").Value = 4170617.2
$ws.Range("J138").Value = 11117333
$ws.Range("L138").Value = 33351999
$ws.Range("N138").Value = -33362279
$ws.Range("H141").Value = 3459.625
$ws.Range("I141").Value = 3097.4285
$ws.Range("J141").Value = 5995
$ws.Range("K141").Value = 9292.2855
$ws.Range("L141").Value = 17985
$ws.Range("M141").Value = -4112.2855
$ws.Range("N141").Value = -28345

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H56").Value = 96661.164
$ws.Range("I56").Value = 69999
$ws.Range("J56").Value = 109992.25
$ws.Range("K56").Value = 69999
$ws.Range("L56").Value = 109992.25
$ws.Range("M56").Value = -69257
$ws.Range("N56").Value = -111476.25
$ws.Range("H61").Value = 8335.037
$ws.Range("I61").Value = 2627.9167
$ws.Range("J61").Value = 12900.733
$ws.Range("K61").Value = 2627.9167
$ws.Range("L61").Value = 12900.733
$ws.Range("M61").Value = -2415.9167
$ws.Range("N61").Value = -13324.733
$ws.Range("H74").Value = 27309.441
$ws.Range("I74").Value = 42419.88
$ws.Range("K74").Value = 42419.88
$ws.Range("M74").Value = -41545.88
$ws.Range("H77").Value = 27309.441
$ws.Range("I77").Value = 42419.88
$ws.Range("K77").Value = 212099.4
$ws.Range("M77").Value = -207731.4
$ws.Range("H132").Value = 986058.4399999999
$ws.Range("I132").Value = 1566472.8
$ws.Range("K132").Value = 4699418.4
$ws.Range("M132").Value = -4696888.4
$ws.Range("H136").Value = 8335.037
$ws.Range("I136").Value = 2627.9167
$ws.Range("J136").Value = 12900.733
$ws.Range("K136").Value = 7883.750100000001
$ws.Range("L136").Value = 38702.199
$ws.Range("M136").Value = -5333.750100000001
$ws.Range("N136").Value = -43802.199

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6802325
$ws.Range("I86").Value = 14718040
$ws.Range("J86").Value = 73967.64999999999
$ws.Range("K86").Value = 14718040
$ws.Range("L86").Value = 73967.64999999999
$ws.Range("M86").Value = -14716917
$ws.Range("N86").Value = -76213.64999999999
$ws.Range("H89").Value = 6802325
$ws.Range("I89").Value = 14718040
$ws.Range("J89").Value = 73967.64999999999
$ws.Range("K89").Value = 73590200
$ws.Range("L89").Value = 369838.25
$ws.Range("M89").Value = -73584584
$ws.Range("N89").Value = -381070.25
$ws.Range("H99").Value = 4135671.2
$ws.Range("I99").Value = 2990.6428
$ws.Range("K99").Value = 2990.6428
$ws.Range("M99").Value = -1492.6428
$ws.Range("H107").Value = 56252820
$ws.Range("I107").Value = 70314216
$ws.Range("K107").Value = 70314216
$ws.Range("M107").Value = -70312296
$ws.Range("H134").Value = 5328.3057
$ws.Range("I134").Value = 1563.0454
$ws.Range("K134").Value = 4689.1362
$ws.Range("M134").Value = -2154.1362

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6088.2554
$ws.Range("I31").Value = 1745.875
$ws.Range("J31").Value = 10619.435
$ws.Range("K31").Value = 1745.875
$ws.Range("L31").Value = 10619.435
$ws.Range("M31").Value = -1450.875
$ws.Range("N31").Value = -11209.435
$ws.Range("H34").Value = 6088.2554
$ws.Range("I34").Value = 1745.875
$ws.Range("J34").Value = 10619.435
$ws.Range("K34").Value = 1745.875
$ws.Range("L34").Value = 10619.435
$ws.Range("M34").Value = -1543.875
$ws.Range("N34").Value = -11023.435
$ws.Range("H58").Value = 10215.762
$ws.Range("I58").Value = 2023
$ws.Range("K58").Value = 2023
$ws.Range("M58").Value = -1820
$ws.Range("H132").Value = 7739.1113
$ws.Range("I132").Value = 5432
$ws.Range("J132").Value = 9584.799999999999
$ws.Range("K132").Value = 16296
$ws.Range("L132").Value = 28754.4
$ws.Range("M132").Value = -13766
$ws.Range("N132").Value = -33814.39999999999
$ws.Range("H134").Value = 5883.6
$ws.Range("I134").Value = 2021.8096
$ws.Range("J134").Value = 11676.286
$ws.Range("K134").Value = 6065.4288
$ws.Range("L134").Value = 35028.858
$ws.Range("M134").Value = -3530.4288
$ws.Range("N134").Value = -40098.858
$ws.Range("H136").Value = 10215.762
$ws.Range("I136").Value = 2023
$ws.Range("K136").Value = 6069
$ws.Range("M136").Value = -3519

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 309
$ws.Range("I7").Value = 283.33334
$ws.Range("K7").Value = 850.0000200000001
$ws.Range("M7").Value = -738.0000200000001
$ws.Range("H122").Value = 708138.6
$ws.Range("I122").Value = 2021202.4
$ws.Range("J122").Value = 1104.3077
$ws.Range("K122").Value = 18190821.6
$ws.Range("L122").Value = 9938.7693
$ws.Range("M122").Value = -18188371.6
$ws.Range("N122").Value = -14838.7693
$ws.Range("H129").Value = 11976898
$ws.Range("I129").Value = 544.3333
$ws.Range("K129").Value = 1632.9999
$ws.Range("M129").Value = 3367.0001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H52").Value = 89600
$ws.Range("J52").Value = 89600
$ws.Range("L52").Value = 89600
$ws.Range("N52").Value = -90118
$ws.Range("H122").Value = 52688476
$ws.Range("I122").Value = 100103120
$ws.Range("J122").Value = 5540.8887
$ws.Range("K122").Value = 300309360
$ws.Range("L122").Value = 16622.6661
$ws.Range("M122").Value = -300306910
$ws.Range("N122").Value = -21522.6661
$ws.Range("H126").Value = 4935.2666
$ws.Range("I126").Value = 2426.4
$ws.Range("K126").Value = 7279.200000000001
$ws.Range("M126").Value = -4809.200000000001
$ws.Range("H132").Value = 3370.625
$ws.Range("I132").Value = 2061.6296
$ws.Range("K132").Value = 6184.888800000001
$ws.Range("M132").Value = -3654.888800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6671662.5
$ws.Range("I61").Value = 14287921
$ws.Range("J61").Value = 7436.875
$ws.Range("K61").Value = 14287921
$ws.Range("L61").Value = 7436.875
$ws.Range("M61").Value = -14287719
$ws.Range("N61").Value = -7840.875
$ws.Range("H96").Value = 51998.5
$ws.Range("J96").Value = 51998.5
$ws.Range("L96").Value = 51998.5
$ws.Range("N96").Value = -57490.5
$ws.Range("H113").Value = 6671662.5
$ws.Range("I113").Value = 14287921
$ws.Range("J113").Value = 7436.875
$ws.Range("K113").Value = 14287921
$ws.Range("L113").Value = 7436.875
$ws.Range("M113").Value = -14285751
$ws.Range("N113").Value = -11776.875
$ws.Range("H122").Value = 9027.521000000001
$ws.Range("I122").Value = 12737.556
$ws.Range("K122").Value = 38212.66800000001
$ws.Range("M122").Value = -35762.66800000001
$ws.Range("H132").Value = 6210.841
$ws.Range("I132").Value = 3744.65
$ws.Range("K132").Value = 11233.95
$ws.Range("M132").Value = -8703.950000000001
$ws.Range("H136").Value = 8099.593
$ws.Range("I136").Value = 5494.4243
$ws.Range("J136").Value = 12193.429
$ws.Range("K136").Value = 16483.2729
$ws.Range("L136").Value = 36580.287
$ws.Range("M136").Value = -13933.2729
$ws.Range("N136").Value = -41680.287
$ws.Range("H141").Value = 68996.625
$ws.Range("J141").Value = 68996.625
$ws.Range("L141").Value = 68996.625
$ws.Range("N141").Value = -79356.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 19097890
$ws.Range("I81").Value = 1251472.9
$ws.Range("K81").Value = 2502945.8
$ws.Range("M81").Value = -2501884.8
$ws.Range("H84").Value = 19097890
$ws.Range("I84").Value = 1251472.9
$ws.Range("K84").Value = 12514729
$ws.Range("M84").Value = -12509425
$ws.Range("H114").Value = 30000
$ws.Range("I114").Value = 30000
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 30000
$ws.Range("L114").Value = 0
$ws.Range("M114").Value = -25661
$ws.Range("N114").ClearContents()
$ws.Range("H122").Value = 71127.13
$ws.Range("I122").Value = 101292.56
$ws.Range("J122").Value = 9288
$ws.Range("K122").Value = 303877.68
$ws.Range("L122").Value = 27864
$ws.Range("M122").Value = -301427.68
$ws.Range("N122").Value = -32764
$ws.Range("H132").Value = 13739.404
$ws.Range("I132").Value = 10558.479
$ws.Range("J132").Value = 17590
$ws.Range("K132").Value = 31675.437
$ws.Range("L132").Value = 52770
$ws.Range("M132").Value = -29145.437
$ws.Range("N132").Value = -57830
$ws.Range("H136").Value = 32761.457
$ws.Range("I136").Value = 1511.5454
$ws.Range("J136").Value = 85645.92
$ws.Range("K136").Value = 4534.6362
$ws.Range("L136").Value = 256937.76
$ws.Range("M136").Value = -1984.6362
$ws.Range("N136").Value = -262037.76
